$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4   = 8.15
    6   = 6.141
    7   = 6.025
    8   = 5.96
    16  = 5.389
    20  = 7.501
    21  = 9.099
    28  = 6.309
    29  = 5.615
    30  = 6.112
    32  = 6.845999999999999
    40  = 9.084
    46  = 6.627
    51  = 5.765000000000001
    52  = 6.116000000000001
    57  = 5.189
    59  = 4.848999999999999
    62  = 5.261
    66  = 5.114999999999999
    73  = 6.703
    74  = 9.013000000000002
    77  = 5.805999999999999
    92  = 5.371
    100 = 6.031000000000001
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 2).Value = $updates[$row]
}
